$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 is the "Umbria" status row. The update (per commit "Aggiornato file
# excel con Umbria") records that the Umbria processing finished:
#  - Ultima rilevazione (B12): bumped one day, 14/03/2018 -> 15/03/2018
#  - Elaborazione (C12): "Non completato. " -> the standard
#    "Completato controlli qualità interni e controlli spaziali" status text
#    used by every other completed region, with the same centered formatting.
#  - Nota (D12) / Fonte Dati (E12): text/content unchanged.

$b12 = $ws.Cells.Item(12, 2)
$b12.NumberFormat = "DD/MM/YY"
$b12.Value = 43174

$c12 = $ws.Cells.Item(12, 3)
$c12.Value = "Completato controlli qualità interni e controlli spaziali"
$c12.HorizontalAlignment = -4108
$c12.VerticalAlignment = -4107
$c12.WrapText = $false

# Move the saved selection/cursor to B13 (matches the saved view state).
$ws.Range("B13").Select() | Out-Null
